$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy formatting from E1 (bold header style) and set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F9 - timestamps as text (must stay plain text, not be parsed as dates/numbers)
$timestamps = @(
    "2021-10-05 10:49:56.144873",
    "2021-10-05 10:49:56.144884",
    "2021-10-05 10:49:56.144888",
    "2021-10-05 10:49:56.144890",
    "2021-10-05 10:49:56.144893",
    "2021-10-05 10:49:56.144896",
    "2021-10-05 10:49:56.144899",
    "2021-10-05 10:49:56.144901"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
